$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = -3
$ws.Range("D4").Value = 45840.80861854634
$ws.Range("E4").Value = -3
$ws.Range("F4").Value = 45840.70297453704
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 45840.80861854191
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 45840.69604166667
$ws.Range("C23").Value = 75
$ws.Range("D23").Value = 45840.80859799153
$ws.Range("E23").Value = 75
$ws.Range("F23").Value = 45840.38353009259
$ws.Range("C33").Value = 2601
$ws.Range("D33").Value = 45840.80859799837
$ws.Range("E33").Value = 2601
$ws.Range("F33").Value = 45840.48958333334
$ws.Range("C56").Value = 193
$ws.Range("D56").Value = 45840.80861854657
$ws.Range("E56").Value = 193
$ws.Range("F56").Value = 45840.70297453704
$ws.Range("C58").Value = 83
$ws.Range("D58").Value = 45840.80861855218
$ws.Range("E58").Value = 83
$ws.Range("F58").Value = 45840.70408564815
$ws.Range("C59").Value = 84
$ws.Range("D59").Value = 45840.80859799736
$ws.Range("E59").Value = 84
$ws.Range("F59").Value = 45840.48582175926
$ws.Range("C63").Value = 105
$ws.Range("D63").Value = 45840.80859799984
$ws.Range("E63").Value = 105
$ws.Range("F63").Value = 45840.51295138889
$ws.Range("C69").Value = 4
$ws.Range("D69").Value = 45840.80861854216
$ws.Range("E69").Value = 4
$ws.Range("F69").Value = 45840.69724537037
$ws.Range("C70").Value = 39
$ws.Range("D70").Value = 45840.80861854237
$ws.Range("E70").Value = 39
$ws.Range("F70").Value = 45840.69724537037
$ws.Range("C81").Value = 166
$ws.Range("D81").Value = 45840.80861854681
$ws.Range("E81").Value = 166
$ws.Range("F81").Value = 45840.70297453704
$ws.Range("C82").Value = 40
$ws.Range("D82").Value = 45840.8085980047
$ws.Range("E82").Value = 40
$ws.Range("F82").Value = 45840.65059027778
$ws.Range("C83").Value = 181
$ws.Range("D83").Value = 45840.80861854704
$ws.Range("E83").Value = 181
$ws.Range("F83").Value = 45840.70297453704
$ws.Range("C85").Value = 48
$ws.Range("D85").Value = 45840.80859800373
$ws.Range("E85").Value = 48
$ws.Range("F85").Value = 45840.62773148148
$ws.Range("C86").Value = 107
$ws.Range("D86").Value = 45840.80861854728
$ws.Range("E86").Value = 107
$ws.Range("F86").Value = 45840.70297453704
$ws.Range("C87").Value = 31
$ws.Range("D87").Value = 45840.80861854751
$ws.Range("E87").Value = 31
$ws.Range("F87").Value = 45840.70297453704
$ws.Range("C90").Value = 18
$ws.Range("D90").Value = 45840.80859800224
$ws.Range("E90").Value = 18
$ws.Range("F90").Value = 45840.56913194444
$ws.Range("C91").Value = 41
$ws.Range("D91").Value = 45840.80859799853
$ws.Range("E91").Value = 41
$ws.Range("F91").Value = 45840.48958333334
$ws.Range("C94").Value = 57
$ws.Range("D94").Value = 45840.80861855553
$ws.Range("E94").Value = 57
$ws.Range("F94").Value = 45840.73303240741
$ws.Range("C96").Value = 81
$ws.Range("D96").Value = 45840.80861855617
$ws.Range("E96").Value = 81
$ws.Range("F96").Value = 45840.75590277778
$ws.Range("C98").Value = 512
$ws.Range("D98").Value = 45840.80861854261
$ws.Range("E98").Value = 512
$ws.Range("F98").Value = 45840.69724537037
$ws.Range("C101").Value = 1043
$ws.Range("D101").Value = 45840.8086185526
$ws.Range("E101").Value = 1043
$ws.Range("F101").Value = 45840.70496527778
$ws.Range("C103").Value = 57
$ws.Range("D103").Value = 45840.80861855284
$ws.Range("E103").Value = 57
$ws.Range("F103").Value = 45840.70496527778
$ws.Range("C106").Value = 92
$ws.Range("D106").Value = 45840.80859800486
$ws.Range("E106").Value = 92
$ws.Range("F106").Value = 45840.65059027778
$ws.Range("C117").Value = 973
$ws.Range("D117").Value = 45840.80859800338
$ws.Range("E117").Value = 973
$ws.Range("F117").Value = 45840.62440972222
$ws.Range("C120").Value = 22
$ws.Range("D120").Value = 45840.80859800502
$ws.Range("E120").Value = 22
$ws.Range("F120").Value = 45840.65059027778
$ws.Range("C123").Value = 349
$ws.Range("D123").Value = 45840.80861854774
$ws.Range("E123").Value = 349
$ws.Range("F123").Value = 45840.70297453704
$ws.Range("C125").Value = 392
$ws.Range("D125").Value = 45840.80859800091
$ws.Range("E125").Value = 392
$ws.Range("F125").Value = 45840.54755787037
$ws.Range("C141").Value = 268
$ws.Range("D141").Value = 45840.80861854798
$ws.Range("E141").Value = 268
$ws.Range("F141").Value = 45840.70297453704
$ws.Range("C145").Value = 37
$ws.Range("D145").Value = 45840.80859799486
$ws.Range("E145").Value = 37
$ws.Range("F145").Value = 45840.42106481481
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 45840.80861855239
$ws.Range("E151").Value = 0
$ws.Range("F151").Value = 45840.70408564815
$ws.Range("C161").Value = 73
$ws.Range("D161").Value = 45840.80859800107
$ws.Range("E161").Value = 73
$ws.Range("F161").Value = 45840.54755787037
$ws.Range("C163").Value = 17
$ws.Range("D163").Value = 45840.8086185346
$ws.Range("E163").Value = 17
$ws.Range("F163").Value = 45840.67962962963
$ws.Range("C175").Value = 103
$ws.Range("D175").Value = 45840.80861855575
$ws.Range("E175").Value = 103
$ws.Range("F175").Value = 45840.7536574074
$ws.Range("C185").Value = 107
$ws.Range("D185").Value = 45840.80859800518
$ws.Range("E185").Value = 107
$ws.Range("F185").Value = 45840.65059027778
$ws.Range("C192").Value = 20
$ws.Range("D192").Value = 45840.80859800665
$ws.Range("E192").Value = 20
$ws.Range("F192").Value = 45840.67293981482
$ws.Range("C195").Value = -6
$ws.Range("D195").Value = 45840.8086185482
$ws.Range("E195").Value = -6
$ws.Range("F195").Value = 45840.70297453704
$ws.Range("C200").Value = 978
$ws.Range("D200").Value = 45840.80859799619
$ws.Range("E200").Value = 978
$ws.Range("F200").Value = 45840.44622685185
$ws.Range("C217").Value = 36
$ws.Range("D217").Value = 45840.80859799434
$ws.Range("E217").Value = 36
$ws.Range("F217").Value = 45840.41077546297
$ws.Range("C218").Value = 46
$ws.Range("D218").Value = 45840.80859800136
$ws.Range("E218").Value = 46
$ws.Range("F218").Value = 45840.54755787037
$ws.Range("C235").Value = 124
$ws.Range("D235").Value = 45840.80861853628
$ws.Range("E235").Value = 124
$ws.Range("F235").Value = 45840.68016203704
$ws.Range("C247").Value = 705
$ws.Range("D247").Value = 45840.80861853909
$ws.Range("E247").Value = 705
$ws.Range("F247").Value = 45840.68561342593
$ws.Range("C255").Value = 1088
$ws.Range("D255").Value = 45840.80861854846
$ws.Range("E255").Value = 1088
$ws.Range("F255").Value = 45840.70297453704
$ws.Range("C270").Value = -2
$ws.Range("D270").Value = 45840.80861854867
$ws.Range("E270").Value = -2
$ws.Range("F270").Value = 45840.70297453704
$ws.Range("C272").Value = 209
$ws.Range("D272").Value = 45840.80861855704
$ws.Range("E272").Value = 209
$ws.Range("F272").Value = 45840.75967592592
$ws.Range("C274").Value = 201
$ws.Range("D274").Value = 45840.8086185412
$ws.Range("E274").Value = 201
$ws.Range("F274").Value = 45840.69284722222
$ws.Range("C283").Value = 126
$ws.Range("D283").Value = 45840.80859800534
$ws.Range("E283").Value = 126
$ws.Range("F283").Value = 45840.65059027778
$ws.Range("C291").Value = 442
$ws.Range("D291").Value = 45840.80861854285
$ws.Range("E291").Value = 442
$ws.Range("F291").Value = 45840.69724537037
$ws.Range("C292").Value = 157
$ws.Range("D292").Value = 45840.80859799637
$ws.Range("E292").Value = 157
$ws.Range("F292").Value = 45840.44622685185
$ws.Range("F295").Value = 45840.70297453704
$ws.Range("C309").Value = 928
$ws.Range("D309").Value = 45840.80859799868
$ws.Range("E309").Value = 928
$ws.Range("F309").Value = 45840.48958333334
$ws.Range("C314").Value = -9
$ws.Range("D314").Value = 45840.80861854419
$ws.Range("E314").Value = -9
$ws.Range("F314").Value = 45840.69724537037
$ws.Range("C315").Value = 73
$ws.Range("D315").Value = 45840.80861853678
$ws.Range("E315").Value = 73
$ws.Range("F315").Value = 45840.68217592593
$ws.Range("C318").Value = 991
$ws.Range("D318").Value = 45840.80859799519
$ws.Range("E318").Value = 991
$ws.Range("F318").Value = 45840.42267361111
$ws.Range("C320").Value = -1
$ws.Range("D320").Value = 45840.80859800679
$ws.Range("E320").Value = -1
$ws.Range("F320").Value = 45840.67293981482
$ws.Range("C321").Value = 61
$ws.Range("D321").Value = 45840.80859800389
$ws.Range("E321").Value = 61
$ws.Range("F321").Value = 45840.62773148148
$ws.Range("C326").Value = 1053
$ws.Range("D326").Value = 45840.80861854913
$ws.Range("E326").Value = 1053
$ws.Range("F326").Value = 45840.70297453704
$ws.Range("C346").Value = 73
$ws.Range("D346").Value = 45840.80861854935
$ws.Range("E346").Value = 73
$ws.Range("F346").Value = 45840.70297453704
$ws.Range("C351").Value = 1082
$ws.Range("D351").Value = 45840.808618537
$ws.Range("E351").Value = 1082
$ws.Range("F351").Value = 45840.68217592593
$ws.Range("C354").Value = -2
$ws.Range("D354").Value = 45840.8085979946
$ws.Range("E354").Value = -2
$ws.Range("F354").Value = 45840.41077546297
$ws.Range("C359").Value = 13
$ws.Range("D359").Value = 45840.80859800289
$ws.Range("E359").Value = 13
$ws.Range("F359").Value = 45840.62070601852
$ws.Range("C363").Value = 405
$ws.Range("D363").Value = 45840.80861855308
$ws.Range("E363").Value = 405
$ws.Range("F363").Value = 45840.70496527778
$ws.Range("C371").Value = 5
$ws.Range("D371").Value = 45840.80859799654
$ws.Range("E371").Value = 5
$ws.Range("F371").Value = 45840.44622685185
$ws.Range("C376").Value = 9
$ws.Range("D376").Value = 45840.80861854956
$ws.Range("E376").Value = 9
$ws.Range("F376").Value = 45840.70297453704
$ws.Range("C379").Value = 0
$ws.Range("D379").Value = 45840.80861854977
$ws.Range("E379").Value = 0
$ws.Range("F379").Value = 45840.70297453704
$ws.Range("C384").Value = 40
$ws.Range("D384").Value = 45840.80861853438
$ws.Range("E384").Value = 40
$ws.Range("F384").Value = 45840.67901620371
$ws.Range("C385").Value = 209
$ws.Range("D385").Value = 45840.8086185365
$ws.Range("E385").Value = 209
$ws.Range("F385").Value = 45840.68038194445
$ws.Range("C390").Value = 142
$ws.Range("D390").Value = 45840.80859800404
$ws.Range("E390").Value = 142
$ws.Range("F390").Value = 45840.62773148148
$ws.Range("C394").Value = 185
$ws.Range("D394").Value = 45840.80861854444
$ws.Range("E394").Value = 185
$ws.Range("F394").Value = 45840.69724537037
$ws.Range("C395").Value = 37
$ws.Range("D395").Value = 45840.80861855196
$ws.Range("E395").Value = 37
$ws.Range("F395").Value = 45840.70306712963
$ws.Range("C402").Value = 100
$ws.Range("D402").Value = 45840.80859799884
$ws.Range("E402").Value = 100
$ws.Range("F402").Value = 45840.48958333334
$ws.Range("C404").Value = 432
$ws.Range("D404").Value = 45840.80861854469
$ws.Range("E404").Value = 432
$ws.Range("F404").Value = 45840.69724537037
$ws.Range("C409").Value = 203
$ws.Range("D409").Value = 45840.80859800156
$ws.Range("E409").Value = 203
$ws.Range("F409").Value = 45840.54755787037
$ws.Range("C410").Value = 1353
$ws.Range("D410").Value = 45840.80861854494
$ws.Range("E410").Value = 1353
$ws.Range("F410").Value = 45840.69724537037
$ws.Range("C418").Value = 707
$ws.Range("D418").Value = 45840.80859800553
$ws.Range("E418").Value = 707
$ws.Range("F418").Value = 45840.65059027778
$ws.Range("C442").Value = 23
$ws.Range("D442").Value = 45840.80861854996
$ws.Range("E442").Value = 23
$ws.Range("F442").Value = 45840.70297453704
$ws.Range("C465").Value = 15
$ws.Range("D465").Value = 45840.80861855639
$ws.Range("E465").Value = 15
$ws.Range("F465").Value = 45840.75590277778
$ws.Range("C469").Value = 2906
$ws.Range("D469").Value = 45840.80861855018
$ws.Range("E469").Value = 2906
$ws.Range("F469").Value = 45840.70297453704
$ws.Range("C472").Value = 14
$ws.Range("D472").Value = 45840.80861855447
$ws.Range("E472").Value = 14
$ws.Range("F472").Value = 45840.7230787037
$ws.Range("C478").Value = -1
$ws.Range("D478").Value = 45840.80861853484
$ws.Range("E478").Value = -1
$ws.Range("F478").Value = 45840.67962962963
$ws.Range("C480").Value = 234
$ws.Range("D480").Value = 45840.80861855725
$ws.Range("E480").Value = 234
$ws.Range("F480").Value = 45840.75967592592
$ws.Range("C485").Value = 8
$ws.Range("D485").Value = 45840.80859800074
$ws.Range("E485").Value = 8
$ws.Range("F485").Value = 45840.53643518518
$ws.Range("C497").Value = 5
$ws.Range("D497").Value = 45840.80861853723
$ws.Range("E497").Value = 5
$ws.Range("F497").Value = 45840.68217592593
$ws.Range("C510").Value = 180
$ws.Range("D510").Value = 45840.80861853746
$ws.Range("E510").Value = 180
$ws.Range("F510").Value = 45840.68217592593
$ws.Range("C528").Value = 130
$ws.Range("D528").Value = 45840.80859799535
$ws.Range("E528").Value = 130
$ws.Range("F528").Value = 45840.42267361111
$ws.Range("C569").Value = 3
$ws.Range("D569").Value = 45840.80861853981
$ws.Range("E569").Value = 3
$ws.Range("F569").Value = 45840.68615740741
$ws.Range("C570").Value = 2496
$ws.Range("D570").Value = 45840.80861853509
$ws.Range("E570").Value = 2496
$ws.Range("F570").Value = 45840.67962962963
$ws.Range("C574").Value = 0
$ws.Range("D574").Value = 45840.80859800648
$ws.Range("E574").Value = 0
$ws.Range("F574").Value = 45840.67111111111
$ws.Range("C581").Value = 19
$ws.Range("D581").Value = 45840.80859800016
$ws.Range("E581").Value = 19
$ws.Range("F581").Value = 45840.51476851852
$ws.Range("C601").Value = -3
$ws.Range("D601").Value = 45840.80861853769
$ws.Range("E601").Value = -3
$ws.Range("F601").Value = 45840.68217592593
$ws.Range("C631").Value = 31
$ws.Range("D631").Value = 45840.80859799013
$ws.Range("E631").Value = 31
$ws.Range("F631").Value = 45840.37972222222
$ws.Range("C657").Value = 1798
$ws.Range("D657").Value = 45840.80859800569
$ws.Range("E657").Value = 1798
$ws.Range("F657").Value = 45840.65059027778
$ws.Range("C658").Value = 4
$ws.Range("D658").Value = 45840.80861854518
$ws.Range("E658").Value = 4
$ws.Range("F658").Value = 45840.69724537037
$ws.Range("C660").Value = 357
$ws.Range("D660").Value = 45840.80861853534
$ws.Range("E660").Value = 357
$ws.Range("F660").Value = 45840.67962962963
$ws.Range("C701").Value = 78
$ws.Range("D701").Value = 45840.80859799901
$ws.Range("E701").Value = 78
$ws.Range("F701").Value = 45840.48958333334
$ws.Range("C720").Value = 397
$ws.Range("D720").Value = 45840.80859799061
$ws.Range("E720").Value = 397
$ws.Range("F720").Value = 45840.37972222222
$ws.Range("C721").Value = -9
$ws.Range("D721").Value = 45840.80859799754
$ws.Range("E721").Value = -9
$ws.Range("F721").Value = 45840.48582175926
$ws.Range("C726").Value = -6
$ws.Range("D726").Value = 45840.80861855745
$ws.Range("E726").Value = -6
$ws.Range("F726").Value = 45840.75967592592
$ws.Range("C732").Value = 170
$ws.Range("D732").Value = 45840.8086185533
$ws.Range("E732").Value = 170
$ws.Range("F732").Value = 45840.70496527778
$ws.Range("C757").Value = 96
$ws.Range("D757").Value = 45840.80861855041
$ws.Range("E757").Value = 96
$ws.Range("F757").Value = 45840.70297453704
$ws.Range("C771").Value = 15
$ws.Range("D771").Value = 45840.8085979918
$ws.Range("E771").Value = 15
$ws.Range("F771").Value = 45840.38353009259
$ws.Range("C778").Value = -4
$ws.Range("D778").Value = 45840.80861853885
$ws.Range("E778").Value = -4
$ws.Range("F778").Value = 45840.68482638889
$ws.Range("C797").Value = 0
$ws.Range("D797").Value = 45840.80861855403
$ws.Range("E797").Value = 0
$ws.Range("F797").Value = 45840.7078587963
$ws.Range("C810").Value = 3
$ws.Range("D810").Value = 45840.80861853389
$ws.Range("E810").Value = 3
$ws.Range("F810").Value = 45840.67560185185
$ws.Range("C811").Value = 5
$ws.Range("D811").Value = 45840.80859800617
$ws.Range("E811").Value = 5
$ws.Range("F811").Value = 45840.66892361111
$ws.Range("C815").Value = 18
$ws.Range("D815").Value = 45840.80861853356
$ws.Range("E815").Value = 18
$ws.Range("F815").Value = 45840.67481481482
$ws.Range("C826").Value = -20
$ws.Range("D826").Value = 45840.80861855469
$ws.Range("E826").Value = -20
$ws.Range("F826").Value = 45840.7236574074
$ws.Range("F866").Value = 45840.7325462963
$ws.Range("C872").Value = 446
$ws.Range("D872").Value = 45840.80859800176
$ws.Range("E872").Value = 446
$ws.Range("F872").Value = 45840.54755787037
$ws.Range("C883").Value = 468
$ws.Range("D883").Value = 45840.80861855511
$ws.Range("E883").Value = 468
$ws.Range("F883").Value = 45840.73126157407
$ws.Range("C884").Value = 49
$ws.Range("D884").Value = 45840.80859799209
$ws.Range("E884").Value = 49
$ws.Range("F884").Value = 45840.38353009259
$ws.Range("C888").Value = 402
$ws.Range("D888").Value = 45840.80861855353
$ws.Range("E888").Value = 402
$ws.Range("F888").Value = 45840.70496527778
$ws.Range("C902").Value = 34
$ws.Range("D902").Value = 45840.8085979967
$ws.Range("E902").Value = 34
$ws.Range("F902").Value = 45840.44622685185
$ws.Range("C909").Value = 3
$ws.Range("D909").Value = 45840.80859799602
$ws.Range("E909").Value = 3
$ws.Range("F909").Value = 45840.44586805555
$ws.Range("C924").Value = 246
$ws.Range("D924").Value = 45840.80859799771
$ws.Range("E924").Value = 246
$ws.Range("F924").Value = 45840.48582175926
$ws.Range("C944").Value = 0
$ws.Range("D944").Value = 45840.80859799346
$ws.Range("E944").Value = 0
$ws.Range("F944").Value = 45840.39509259259
$ws.Range("C962").Value = 320
$ws.Range("D962").Value = 45840.80861855061
$ws.Range("E962").Value = 320
$ws.Range("F962").Value = 45840.70297453704
$ws.Range("C963").Value = 1442
$ws.Range("D963").Value = 45840.80859799918
$ws.Range("E963").Value = 1442
$ws.Range("F963").Value = 45840.48958333334
$ws.Range("C996").Value = 16
$ws.Range("D996").Value = 45840.80861854541
$ws.Range("E996").Value = 16
$ws.Range("F996").Value = 45840.69724537037
$ws.Range("C1002").Value = 60
$ws.Range("D1002").Value = 45840.80861855081
$ws.Range("E1002").Value = 60
$ws.Range("F1002").Value = 45840.70297453704
$ws.Range("C1017").Value = 519
$ws.Range("D1017").Value = 45840.80859800454
$ws.Range("E1017").Value = 519
$ws.Range("F1017").Value = 45840.6453587963
$ws.Range("C1024").Value = 84
$ws.Range("D1024").Value = 45840.808618551
$ws.Range("E1024").Value = 84
$ws.Range("F1024").Value = 45840.70297453704
$ws.Range("C1062").Value = 164
$ws.Range("D1062").Value = 45840.80859800585
$ws.Range("E1062").Value = 164
$ws.Range("F1062").Value = 45840.65059027778
$ws.Range("C1109").Value = 3
$ws.Range("D1109").Value = 45840.80861855682
$ws.Range("E1109").Value = 3
$ws.Range("F1109").Value = 45840.67798611111
$ws.Range("C1110").Value = 99
$ws.Range("D1110").Value = 45840.80859800192
$ws.Range("E1110").Value = 99
$ws.Range("F1110").Value = 45840.54755787037
$ws.Range("C1122").Value = 7
$ws.Range("D1122").Value = 45840.8085980042
$ws.Range("E1122").Value = 7
$ws.Range("F1122").Value = 45840.62773148148
$ws.Range("C1126").Value = 612
$ws.Range("D1126").Value = 45840.80861854563
$ws.Range("E1126").Value = 612
$ws.Range("F1126").Value = 45840.69724537037
$ws.Range("C1152").Value = 1
$ws.Range("D1152").Value = 45840.80859800323
$ws.Range("E1152").Value = 1
$ws.Range("F1152").Value = 45840.62177083334
$ws.Range("C1154").Value = 30
$ws.Range("D1154").Value = 45840.80859799787
$ws.Range("E1154").Value = 30
$ws.Range("F1154").Value = 45840.48582175926
$ws.Range("C1155").Value = -13
$ws.Range("D1155").Value = 45840.80859799803
$ws.Range("E1155").Value = -13
$ws.Range("F1155").Value = 45840.48582175926
$ws.Range("C1193").Value = -3
$ws.Range("D1193").Value = 45840.80861855121
$ws.Range("E1193").Value = -3
$ws.Range("F1193").Value = 45840.70297453704
$ws.Range("C1195").Value = 8
$ws.Range("D1195").Value = 45840.80859799503
$ws.Range("E1195").Value = 8
$ws.Range("F1195").Value = 45840.42232638889
$ws.Range("C1253").Value = 1095
$ws.Range("D1253").Value = 45840.80859799821
$ws.Range("E1253").Value = 1095
$ws.Range("F1253").Value = 45840.48582175926
$ws.Range("C1255").Value = 7
$ws.Range("D1255").Value = 45840.80859799999
$ws.Range("E1255").Value = 7
$ws.Range("F1255").Value = 45840.51420138889
$ws.Range("C1285").Value = 245
$ws.Range("D1285").Value = 45840.80859799966
$ws.Range("E1285").Value = 245
$ws.Range("F1285").Value = 45840.51243055556
$ws.Range("C1332").Value = 33
$ws.Range("D1332").Value = 45840.80861853933
$ws.Range("E1332").Value = 33
$ws.Range("F1332").Value = 45840.68561342593
$ws.Range("C1342").Value = 1133
$ws.Range("D1342").Value = 45840.80859800601
$ws.Range("E1342").Value = 1133
$ws.Range("F1342").Value = 45840.65059027778
$ws.Range("C1369").Value = 43
$ws.Range("D1369").Value = 45840.80861854611
$ws.Range("E1369").Value = 43
$ws.Range("F1369").Value = 45840.69929398148
$ws.Range("C1388").Value = 221
$ws.Range("D1388").Value = 45840.80859799376
$ws.Range("E1388").Value = 221
$ws.Range("F1388").Value = 45840.39947916667
$ws.Range("C1401").Value = 0
$ws.Range("D1401").Value = 45840.80859800257
$ws.Range("E1401").Value = 0
$ws.Range("F1401").Value = 45840.39284722223
$ws.Range("C1403").Value = -4
$ws.Range("D1403").Value = 45840.80861855425
$ws.Range("E1403").Value = -4
$ws.Range("F1403").Value = 45840.72278935185
$ws.Range("C1411").Value = 67
$ws.Range("D1411").Value = 45840.80859800633
$ws.Range("E1411").Value = 67
$ws.Range("F1411").Value = 45840.66932870371
$ws.Range("C1412").Value = -13
$ws.Range("D1412").Value = 45840.80861853957
$ws.Range("E1412").Value = -13
$ws.Range("F1412").Value = 45840.68561342593
$ws.Range("C1421").Value = 49
$ws.Range("D1421").Value = 45840.808598007
$ws.Range("E1421").Value = 49
$ws.Range("F1421").Value = 45840.67443287037
$ws.Range("C1446").Value = 27
$ws.Range("D1446").Value = 45840.80859799092
$ws.Range("E1446").Value = 27
$ws.Range("F1446").Value = 45840.37972222222
$ws.Range("C1448").Value = -3
$ws.Range("D1448").Value = 45840.80859799934
$ws.Range("E1448").Value = -3
$ws.Range("F1448").Value = 45840.48958333334
$ws.Range("C1451").Value = 3
$ws.Range("D1451").Value = 45840.80861853413
$ws.Range("E1451").Value = 3
$ws.Range("F1451").Value = 45840.67635416667
$ws.Range("C1483").Value = 33
$ws.Range("D1483").Value = 45840.80859799687
$ws.Range("E1483").Value = 33
$ws.Range("F1483").Value = 45840.44622685185
$ws.Range("C1502").Value = 2
$ws.Range("D1502").Value = 45840.80861855661
$ws.Range("E1502").Value = 2
$ws.Range("F1502").Value = 45840.37972222222
$ws.Range("C1507").Value = 188
$ws.Range("D1507").Value = 45840.80861855143
$ws.Range("E1507").Value = 188
$ws.Range("F1507").Value = 45840.70297453704
$ws.Range("C1558").Value = 100
$ws.Range("D1558").Value = 45840.80859800436
$ws.Range("E1558").Value = 100
$ws.Range("F1558").Value = 45840.62773148148
$ws.Range("C1567").Value = -140
$ws.Range("D1567").Value = 45840.80861853558
$ws.Range("E1567").Value = -140
$ws.Range("F1567").Value = 45840.67962962963
$ws.Range("C1576").Value = 24
$ws.Range("D1576").Value = 45840.80859799566
$ws.Range("E1576").Value = 24
$ws.Range("F1576").Value = 45840.42353009259
$ws.Range("C1594").Value = 53
$ws.Range("D1594").Value = 45840.80861855175
$ws.Range("E1594").Value = 53
$ws.Range("F1594").Value = 45840.70297453704
$ws.Range("C1597").Value = 4190
$ws.Range("D1597").Value = 45840.80861855596
$ws.Range("E1597").Value = 4190
$ws.Range("F1597").Value = 45840.7536574074
$ws.Range("C1635").Value = 1734
$ws.Range("D1635").Value = 45840.80861853581
$ws.Range("E1635").Value = 1734
$ws.Range("F1635").Value = 45840.67962962963
$ws.Range("C1655").Value = -3
$ws.Range("D1655").Value = 45840.80861853791
$ws.Range("E1655").Value = -3
$ws.Range("F1655").Value = 45840.68217592593
$ws.Range("C1656").Value = 74
$ws.Range("D1656").Value = 45840.80859800208
$ws.Range("E1656").Value = 74
$ws.Range("F1656").Value = 45840.54755787037
$ws.Range("C1737").Value = 6
$ws.Range("D1737").Value = 45840.80861854168
$ws.Range("E1737").Value = 6
$ws.Range("F1737").Value = 45840.69471064815
$ws.Range("C1739").Value = 17
$ws.Range("D1739").Value = 45840.80859800307
$ws.Range("E1739").Value = 17
$ws.Range("F1739").Value = 45840.62070601852
$ws.Range("C1753").Value = 0
$ws.Range("D1753").Value = 45840.80861855491
$ws.Range("E1753").Value = 0
$ws.Range("F1753").Value = 45840.72748842592
$ws.Range("C1786").Value = 0
$ws.Range("D1786").Value = 45840.80861854142
$ws.Range("E1786").Value = 0
$ws.Range("F1786").Value = 45840.69409722222
$ws.Range("C1823").Value = 29
$ws.Range("D1823").Value = 45840.8085979955
$ws.Range("E1823").Value = 29
$ws.Range("F1823").Value = 45840.42267361111
$ws.Range("C1844").Value = 399
$ws.Range("D1844").Value = 45840.8085979995
$ws.Range("E1844").Value = 399
$ws.Range("F1844").Value = 45840.48958333334
$ws.Range("C1925").Value = 23
$ws.Range("D1925").Value = 45840.80859799319
$ws.Range("E1925").Value = 23
$ws.Range("F1925").Value = 45840.39346064815
$ws.Range("C1947").Value = 2
$ws.Range("D1947").Value = 45840.80859799405
$ws.Range("E1947").Value = 2
$ws.Range("F1947").Value = 45840.40810185186
$ws.Range("C1948").Value = 1
$ws.Range("D1948").Value = 45840.80861853814
$ws.Range("E1948").Value = 1
$ws.Range("F1948").Value = 45840.68217592593
$ws.Range("C2029").Value = 2
$ws.Range("D2029").Value = 45840.80859800032
$ws.Range("E2029").Value = 2
$ws.Range("F2029").Value = 45840.51508101852
$ws.Range("C2067").Value = 49
$ws.Range("D2067").Value = 45840.80861854587
$ws.Range("E2067").Value = 49
$ws.Range("F2067").Value = 45840.69724537037
$ws.Range("C2102").Value = -3
$ws.Range("D2102").Value = 45840.80861853837
$ws.Range("E2102").Value = -3
$ws.Range("F2102").Value = 45840.68217592593
$ws.Range("C2113").Value = 0
$ws.Range("D2113").Value = 45840.80861853859
$ws.Range("E2113").Value = 0
$ws.Range("F2113").Value = 45840.68248842593
$ws.Range("C2239").Value = 44
$ws.Range("D2239").Value = 45840.8086358626
$ws.Range("E2239").Value = 44
$ws.Range("F2239").Value = 45840.75967592592
$ws.Range("C2327").Value = 5
$ws.Range("D2327").Value = 45840.80861854018
$ws.Range("E2327").Value = 5
$ws.Range("F2327").Value = 45840.68648148148
$ws.Range("C2328").Value = 16
$ws.Range("D2328").Value = 45840.80861854043
$ws.Range("E2328").Value = 16
$ws.Range("F2328").Value = 45840.68708333333
$ws.Range("C2386").Value = 0
$ws.Range("D2386").Value = 45840.80859800058
$ws.Range("E2386").Value = 0
$ws.Range("F2386").Value = 45840.53359953704
$ws.Range("C2416").Value = 48
$ws.Range("D2416").Value = 45840.8085980024
$ws.Range("E2416").Value = 48
$ws.Range("F2416").Value = 45840.57128472222
$ws.Range("C2420").Value = 0
$ws.Range("D2420").Value = 45840.80859799237
$ws.Range("E2420").Value = 0
$ws.Range("F2420").Value = 45840.38622685185
$ws.Range("C2465").Value = -13
$ws.Range("D2465").Value = 45840.80859800356
$ws.Range("E2465").Value = -13
$ws.Range("F2465").Value = 45840.62440972222
$ws.Range("C2474").Value = 12
$ws.Range("D2474").Value = 45840.80861854095
$ws.Range("E2474").Value = 12
$ws.Range("F2474").Value = 45840.68773148148
$ws.Range("C2483").Value = 841
$ws.Range("D2483").Value = 45840.80861855375
$ws.Range("E2483").Value = 841
$ws.Range("F2483").Value = 45840.70496527778
$ws.Range("C2487").Value = 4
$ws.Range("D2487").Value = 45840.80859799265
$ws.Range("E2487").Value = 4
$ws.Range("F2487").Value = 45840.39229166666
$ws.Range("C2537").Value = 0
$ws.Range("D2537").Value = 45840.80859799123
$ws.Range("E2537").Value = 0
$ws.Range("F2537").Value = 45840.38086805555
$ws.Range("C2541").Value = 11
$ws.Range("D2541").Value = 45840.80861853605
$ws.Range("E2541").Value = 11
$ws.Range("F2541").Value = 45840.68011574074
$ws.Range("C2592").Value = 0
$ws.Range("D2592").Value = 45840.80859799703
$ws.Range("E2592").Value = 0
$ws.Range("F2592").Value = 45840.46775462963
$ws.Range("C2593").Value = 0
$ws.Range("D2593").Value = 45840.8085979972
$ws.Range("E2593").Value = 0
$ws.Range("F2593").Value = 45840.46829861111

# New rows 2600-2602
$ws.Range("A2600").Value = 43993277
$ws.Range("B2600").Value = 1
$ws.Range("C2600").Value = 0
$ws.Range("D2600").Value = 45840.8085979957
$ws.Range("D2600").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2600").Value = 0
$ws.Range("G2600").Value = 0
$ws.Range("H2600").Value = "Consistente"

$ws.Range("A2601").Value = 43996292
$ws.Range("B2601").Value = 1
$ws.Range("C2601").Value = 0
$ws.Range("D2601").Value = 45840.80861853656
$ws.Range("D2601").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2601").Value = 0
$ws.Range("F2601").Value = 45840.68118055556
$ws.Range("F2601").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G2601").Value = 0
$ws.Range("H2601").Value = "Consistente"

$ws.Range("A2602").Value = 44002959
$ws.Range("B2602").Value = 1
$ws.Range("C2602").Value = -3
$ws.Range("D2602").Value = 45840.8086185538
$ws.Range("D2602").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2602").Value = -3
$ws.Range("F2602").Value = 45840.70496527778
$ws.Range("F2602").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G2602").Value = 0
$ws.Range("H2602").Value = "Consistente"

Write-Host "All changes applied"
